$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank row 2 so that the data rows (previously 3:16) shift up to (2:15),
# turning the gap between the header and the sorted data into a single contiguous table.
$ws.Rows("2").Delete()

# Re-apply the sort on the data (now A2:B14) so the worksheet's recorded sort
# range/condition follow the shifted rows instead of the stale A3:B15 range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A14"))
$ws.Sort.SetRange($ws.Range("A2:B14"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Update the selection to the header row, matching the post-edit selection state.
$ws.Range("A1:B1").Select()
